# Auto-generated edit script for LOM3015.xlsx content restructuring
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LabelCell($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.Font.Bold = $true
    $r.WrapText = $false
}

function Set-BCell($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.Font.Bold = $false
    $r.WrapText = $true
}

function Set-CCell($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.Font.Bold = $false
    $r.WrapText = $true
    $r.Font.Color = 255
}

# Clear the region that is being restructured (rows 10-26) before rebuilding it
$ws.Range("A10:C26").Clear()

# Row 10
Set-LabelCell 'A10' 'Objetivos:'
Set-BCell 'B10' 'Esta disciplina visa apresentar fundamentos de termodinâmica aplicada à área de ciência e engenharia de materiais. Especial ênfase é dada à energia na forma e calor para aquecimento de sistemas termodinâmicos; cálculos de variação de entalpia; entropia e energia de Gibbs de elementos e compostos em mudanças de estado; cálculos de variação de entalpia; entropia e energia de Gibbs de reação; aplicação da propriedade  energia de Gibbs para avaliação de transformações espontâneas e em equilíbrio; fundamentos de termodinâmica de soluções; cálculos de condições de equilíbrio em sistemas heterogêneos. Apresenta-se também as principais diferenças entre esta disciplina e a disciplina de Termodinâmica de Máquinas.'
Set-CCell 'C10' 'Esta disciplina visa apresentar fundamentos de termodinâmica aplicada à área de ciência e engenharia de materiais. Especial ênfase é dada à energia na forma e calor para aquecimento de sistemas termodinâmicos; cálculos de variação de entalpia; entropia e energia de Gibbs de elementos e compostos em mudanças de estado; cálculos de variação de entalpia; entropia e energia de Gibbs de reação; aplicação da propriedade  energia de Gibbs para avaliação de transformações espontâneas e em equilíbrio; fundamentos de termodinâmica de soluções; cálculos de condições de equilíbrio em sistemas heterogêneos. Apresenta-se também as principais diferenças entre esta disciplina e a disciplina de Termodinâmica de Máquinas.'
$ws.Rows.Item(10).RowHeight = 60

# Row 11
Set-LabelCell 'A11' 'Objectives:'
$ws.Rows.Item(11).RowHeight = 60

# Row 12
Set-LabelCell 'A12' 'Docentes responsáveis:'
$ws.Rows.Item(12).AutoFit()

# Row 13
Set-BCell 'B13' '3577649 - Carlos Angelo Nunes'
Set-CCell 'C13' '3577649 - Carlos Angelo Nunes'
$ws.Rows.Item(13).AutoFit()

# Row 14
Set-BCell 'B14' '1176388 - Luiz Tadeu Fernandes Eleno'
Set-CCell 'C14' '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Rows.Item(14).AutoFit()

# Row 15
Set-LabelCell 'A15' 'Programa resumido:'
Set-BCell 'B15' '1) Introdução; 2) 1a Lei da Termodinâmica 3) 2a e 3a Leis da Termodinâmica 4) Equilíbrio heterogêneo: composição variável da fase gasosa; 5) Equilíbrio heterogêneo: composição variável da fase condensada;'
Set-CCell 'C15' '1) Introdução; 2) 1a Lei da Termodinâmica 3) 2a e 3a Leis da Termodinâmica 4) Equilíbrio heterogêneo: composição variável da fase gasosa; 5) Equilíbrio heterogêneo: composição variável da fase condensada;'
$ws.Rows.Item(15).RowHeight = 60

# Row 16
Set-LabelCell 'A16' 'Short syllabus:'
$ws.Rows.Item(16).RowHeight = 60

# Row 17
Set-LabelCell 'A17' 'Programa:'
Set-BCell 'B17' '1- Introdução: sistema; vizinhanças; fases; equilíbrio; fronteiras adiabáticas e diatérmicas; processos reversíveis e irreversíveis; estado termodinâmico; mudança de estado; processos cíclicos; equação de estado; calor; trabalho.2- A 1ª lei de Termodinâmica: energia interna; capacidades térmicas; entalpia; entalpia de transformação de fases; entalpia de formação e de reação; entalpia de reação em função da temperatura (introdução ao loop termodinâmico).3- A 2ª e 3ª leis da Termodinâmica: Dispersão de energia e entropia; entropia no zero absoluto; entropia de reação; entropia de reação em função da temperatura; desigualdade de Clausius; critérios de espontaneidade e equilíbrio; energia de Gibbs; energia de Helmholtz; energia de Gibbs de reação em função da temperatura; equação de Gibbs-Helmholtz.4- Equilíbrio heterogêneo: composição variável da fase gasosa: mistura de gases ideais; lei de Dalton; energia de Gibbs de um gás ideal; pressão de equilíbrio em sistemas metal-óxido-O2(g).5- Equilíbrio heterogêneo: composição variável da fase condensada: fugacidade; atividade termodinâmica; soluções e grandezas parciais molares; potencial químico; modelos de soluções; propriedades termodinâmicas de excesso'
Set-CCell 'C17' '1- Introdução: sistema; vizinhanças; fases; equilíbrio; fronteiras adiabáticas e diatérmicas; processos reversíveis e irreversíveis; estado termodinâmico; mudança de estado; processos cíclicos; equação de estado; calor; trabalho.2- A 1ª lei de Termodinâmica: energia interna; capacidades térmicas; entalpia; entalpia de transformação de fases; entalpia de formação e de reação; entalpia de reação em função da temperatura (introdução ao loop termodinâmico).3- A 2ª e 3ª leis da Termodinâmica: Dispersão de energia e entropia; entropia no zero absoluto; entropia de reação; entropia de reação em função da temperatura; desigualdade de Clausius; critérios de espontaneidade e equilíbrio; energia de Gibbs; energia de Helmholtz; energia de Gibbs de reação em função da temperatura; equação de Gibbs-Helmholtz.4- Equilíbrio heterogêneo: composição variável da fase gasosa: mistura de gases ideais; lei de Dalton; energia de Gibbs de um gás ideal; pressão de equilíbrio em sistemas metal-óxido-O2(g).5- Equilíbrio heterogêneo: composição variável da fase condensada: fugacidade; atividade termodinâmica; soluções e grandezas parciais molares; potencial químico; modelos de soluções; propriedades termodinâmicas de excesso'
$ws.Rows.Item(17).RowHeight = 120

# Row 18
Set-LabelCell 'A18' 'Syllabus:'
$ws.Rows.Item(18).RowHeight = 120

# Row 19
Set-LabelCell 'A19' 'Avaliação:'
$ws.Rows.Item(19).AutoFit()

# Row 20
Set-LabelCell 'A20' 'Método:'
Set-BCell 'B20' 'Esta é uma disciplina fundamental, exigindo dedicação individual para assimilação de definições e conceitos. Isto envolve leitura concentrada e realização de exercícios numéricos.'
Set-CCell 'C20' 'Esta é uma disciplina fundamental, exigindo dedicação individual para assimilação de definições e conceitos. Isto envolve leitura concentrada e realização de exercícios numéricos.'
$ws.Rows.Item(20).RowHeight = 60

# Row 21
Set-LabelCell 'A21' 'Critério:'
Set-BCell 'B21' 'Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão:NF=(P1+2*P2)/3'
Set-CCell 'C21' 'Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão:NF=(P1+2*P2)/3'
$ws.Rows.Item(21).RowHeight = 60

# Row 22
Set-LabelCell 'A22' 'Norma de recuperação:'
Set-BCell 'B22' 'Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'
Set-CCell 'C22' 'Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'
$ws.Rows.Item(22).RowHeight = 60

# Row 23
Set-LabelCell 'A23' 'Bibliografia:'
Set-BCell 'B23' '1) Johnson, D.L. & Stracher, G.B., Thermodynamic - Loop Applications in Materials Systems, vols.1 e 2, The Minerals, Metals & Materials Society, 1995. ISBN 0-87339-270-1.2) P. Atkins & J. de Paula. Físico-Química, Livros Técnicos e Científicos Editora S.A., 2008. ISBN 978-85-216-1600-9.3) S.Stolen, T.Grande. Chemical Thermodynamics of Materials, John Wiley & Sons, Ltd. 2005. ISBN 978-0-471-49230-6.4) R. DeHoff. Thermodynamics in Materials Science. Taylor & Francis Group, 2006. ISBN 978-0-8493-4065-9.5) Y.A. Chang & W.A. Oates. Materials Thermodynamics, John Wiley & Sons, 2010. ISBN 978-0-470-48414-2.'
Set-CCell 'C23' '1) Johnson, D.L. & Stracher, G.B., Thermodynamic - Loop Applications in Materials Systems, vols.1 e 2, The Minerals, Metals & Materials Society, 1995. ISBN 0-87339-270-1.2) P. Atkins & J. de Paula. Físico-Química, Livros Técnicos e Científicos Editora S.A., 2008. ISBN 978-85-216-1600-9.3) S.Stolen, T.Grande. Chemical Thermodynamics of Materials, John Wiley & Sons, Ltd. 2005. ISBN 978-0-471-49230-6.4) R. DeHoff. Thermodynamics in Materials Science. Taylor & Francis Group, 2006. ISBN 978-0-8493-4065-9.5) Y.A. Chang & W.A. Oates. Materials Thermodynamics, John Wiley & Sons, 2010. ISBN 978-0-470-48414-2.'
$ws.Rows.Item(23).RowHeight = 120

# Row 24
Set-LabelCell 'A24' 'Requisitos:'
$ws.Rows.Item(24).AutoFit()

# Row 25
Set-BCell 'B25' 'LOB1004 -  Cálculo II  (Requisito fraco)
'
Set-CCell 'C25' 'LOB1004 -  Cálculo II  (Requisito fraco)
'
$ws.Rows.Item(25).RowHeight = 30

# Row 26
Set-BCell 'B26' 'LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)
'
Set-CCell 'C26' 'LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)
'
$ws.Rows.Item(26).RowHeight = 30

# Fix column B width (pre-existing overlapping col definition bug)
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(3).ColumnWidth
